{"js": "// Word template edit: turn the hard-coded DIPA year \"2025\" into the\n// \"{{ tahun_anggaran }}\" Jinja-style placeholder in both places the\n// DIPA number \"DIPA-26.13.2.452558/...\" appears in the report body.\n\nconst body = context.document.body;\n\n// --- Location 1 -----------------------------------------------------\n// \" Nomor DIPA-26.13.2.452558/2025, dan Petunjuk Operasional Kegiatan\n// (POK) Balai Besar Pelatihan Vokasi dan Produktivitas Bekasi Tahun\n// Anggaran \" -> the literal \"2025\" becomes \"{{ tahun_anggaran }}\".\nconst loc1 = body.search(\"DIPA-26.13.2.452558/2025,\", { matchCase: true });\nloc1.load(\"items\");\nawait context.sync();\n\nif (loc1.items.length === 0) {\n  throw new Error(\"Location 1 text (DIPA-26.13.2.452558/2025,) not found\");\n}\nloc1.items[0].insertText(\n  \"DIPA-26.13.2.452558/{{ tahun_anggaran }},\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- Location 2 -----------------------------------------------------\n// \" }} Nomor : DIPA-26.13.2.452558/ 2025 tanggal 2 Desember 2024.\" ->\n// the trailing \" 2025 tanggal 2 Desember 2024\" is replaced by the same\n// placeholder, keeping the final period.\nconst loc2 = body.search(\n  \"DIPA-26.13.2.452558/ 2025 tanggal 2 Desember 2024.\",\n  { matchCase: true }\n);\nloc2.load(\"items\");\nawait context.sync();\n\nif (loc2.items.length === 0) {\n  throw new Error(\n    \"Location 2 text (DIPA-26.13.2.452558/ 2025 tanggal 2 Desember 2024.) not found\"\n  );\n}\nloc2.items[0].insertText(\n  \"DIPA-26.13.2.452558/{{ tahun_anggaran }}.\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Word template edit: turn the hard-coded DIPA year \"2025\" into the\n# \"{{ tahun_anggaran }}\" Jinja-style placeholder in both places the\n# DIPA number \"DIPA-26.13.2.452558/...\" appears in the report body.\n\n$d = $word.ActiveDocument\n\n# --- Location 1 -------------------------------------------------------\n# \" Nomor DIPA-26.13.2.452558/2025, dan Petunjuk Operasional Kegiatan\n# (POK) Balai Besar Pelatihan Vokasi dan Produktivitas Bekasi Tahun\n# Anggaran \" -> the literal \"2025\" becomes \"{{ tahun_anggaran }}\".\n$rng1 = $d.Content\n$found1 = $rng1.Find.Execute(\"DIPA-26.13.2.452558/2025,\")\nif ($found1) {\n    $rng1.Text = \"DIPA-26.13.2.452558/{{ tahun_anggaran }},\"\n} else {\n    throw \"Location 1 text (DIPA-26.13.2.452558/2025,) not found\"\n}\n\n# --- Location 2 -------------------------------------------------------\n# \" }} Nomor : DIPA-26.13.2.452558/ 2025 tanggal 2 Desember 2024.\" ->\n# the trailing \" 2025 tanggal 2 Desember 2024\" is replaced by the same\n# placeholder, keeping the final period.\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\"DIPA-26.13.2.452558/ 2025 tanggal 2 Desember 2024.\")\nif ($found2) {\n    $rng2.Text = \"DIPA-26.13.2.452558/{{ tahun_anggaran }}.\"\n} else {\n    throw \"Location 2 text (DIPA-26.13.2.452558/ 2025 tanggal 2 Desember 2024.) not found\"\n}\n"}
